# This script turns the 10-row "ランサーズ" listing sheet into a 2-row
# sheet: the two oldest postings (previously rows 3 and 9) are re-scraped
# with a refreshed timestamp and become the new rows 2 and 3, while every
# other posting (previously rows 4-11) is dropped from the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2026-02-15 06:48:38"

# New row 2 <- content that used to live in row 3 (Salesforce consulting ad)
$ws.Range("A2").Value = $newTimestamp
$ws.Range("B2").Value = "【急募】Salesforce 認定 Service Cloud コンサルタントを探しています!"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5492003"
$ws.Range("G2").Value = 43
$ws.Range("H2").Value = "◆コンサル"

# New row 3 <- content that used to live in row 9 (Ticketmaster trouble ad);
# this posting has no "skill summary" entry, so H3 must stay empty.
$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "【緊急】Ticketmasterのチケット購入後のトラブル解決依頼"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5491983"
$ws.Range("G3").Value = 10
$ws.Range("H3").ClearContents()

# Drop the remaining old listings (previously rows 4-11); this also shrinks
# the sheet's used range/dimension down to A1:H3 automatically.
$ws.Rows("4:11").Delete()

# Narrow the "価格" column slightly.
$ws.Columns("D").ColumnWidth = 25.2

# Rebuild the hyperlinks collection so it only contains the two links that
# belong to the surviving rows, pointing at their (possibly new) URLs.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5492003")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5491983")
